$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    if ($val -match '^-?\d+(\.\d+)?$') {
        $ws.Range($addr).Value = "'" + $val
    } else {
        $ws.Range($addr).Value = $val
    }
}

Set-TextValue $ws "D2" "72.225.40"
Set-TextValue $ws "E2" "  +4.56%  "

Set-TextValue $ws "D3" "4.031.76"
Set-TextValue $ws "E3" "  +4.21%  "

Set-TextValue $ws "E4" "  +0.09%  "

Set-TextValue $ws "D5" "538.05"
Set-TextValue $ws "E5" "  +3.53%  "

Set-TextValue $ws "D6" "153.60"
Set-TextValue $ws "E6" "  +8.10%  "

Set-TextValue $ws "D7" "0.696"
Set-TextValue $ws "E7" "  +14.40%  "

Set-TextValue $ws "E8" "  +0.06%  "

Set-TextValue $ws "D9" "0.755"
Set-TextValue $ws "E9" "  +4.99%  "

Set-TextValue $ws "E10" "  +2.55%  "

Set-TextValue $ws "D11" "0.0000328"
Set-TextValue $ws "E11" "  +0.59%  "

Set-TextValue $ws "D12" "47.89"
Set-TextValue $ws "E12" "  +14.56%  "

Set-TextValue $ws "D13" "10.80"
Set-TextValue $ws "E13" "  +4.08%  "

Set-TextValue $ws "D14" "4.662.28"
Set-TextValue $ws "E14" "  +4.47%  "

Set-TextValue $ws "D15" "4.032.27"
Set-TextValue $ws "E15" "  +3.94%  "

Set-TextValue $ws "D16" "14.23"
Set-TextValue $ws "E16" "  +1.46%  "

Set-TextValue $ws "D17" "20.62"
Set-TextValue $ws "E17" "  -2.79%  "

Set-TextValue $ws "E18" "  -0.40%  "

Set-TextValue $ws "D19" "1.20"
Set-TextValue $ws "E19" "  -0.05%  "

Set-TextValue $ws "D20" "71.985.80"
Set-TextValue $ws "E20" "  +4.52%  "

Set-TextValue $ws "D21" "433.72"
Set-TextValue $ws "E21" "  +3.95%  "

Set-TextValue $ws "D22" "99.04"
Set-TextValue $ws "E22" "  +13.77%  "

Set-TextValue $ws "D23" "3.58"
Set-TextValue $ws "E23" "  +2.95%  "

Set-TextValue $ws "B24" "PancakeSwap"
Set-TextValue $ws "C24" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws "D24" "4.23"
Set-TextValue $ws "E24" "  +5.27%  "

Set-TextValue $ws "B25" "InternetComputer(DFINITY)"
Set-TextValue $ws "C25" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D25" "14.56"
Set-TextValue $ws "E25" "  +3.86%  "

Set-TextValue $ws "D26" "11.16"
Set-TextValue $ws "E26" "  -6.60%  "

Set-TextValue $ws "D27" "10.87"
Set-TextValue $ws "E27" "  +3.77%  "

Set-TextValue $ws "E28" "  +30.46%  "

Set-TextValue $ws "E29" "  +2.96%  "

Set-TextValue $ws "D30" "37.04"
Set-TextValue $ws "E30" "  +3.97%  "

Set-TextValue $ws "D31" "13.49"
Set-TextValue $ws "E31" "  +0.82%  "

Set-TextValue $ws "D32" "0.131"
Set-TextValue $ws "E32" "  +5.12%  "

Set-TextValue $ws "D33" "684.97"
Set-TextValue $ws "E33" "  +1.02%  "

Set-TextValue $ws "D34" "6.84"
Set-TextValue $ws "E34" "  -2.16%  "

Set-TextValue $ws "D35" "66.61"
Set-TextValue $ws "E35" "  +0.27%  "

Set-TextValue $ws "D36" "42.51"
Set-TextValue $ws "E36" "  +7.95%  "

Set-TextValue $ws "D37" "0.428"
Set-TextValue $ws "E37" "  -3.65%  "

Set-TextValue $ws "B38" "Kaspa"
Set-TextValue $ws "C38" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D38" "0.156"
Set-TextValue $ws "E38" "  +6.41%  "

Set-TextValue $ws "B39" "WEMIXToken"
Set-TextValue $ws "C39" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws "D39" "3.54"
Set-TextValue $ws "E39" "  +12.78%  "

Set-TextValue $ws "D40" "0.0₃0834"
Set-TextValue $ws "E40" "  -3.21%  "

Set-TextValue $ws "D41" "3.47"
Set-TextValue $ws "E41" "  +2.29%  "

Set-TextValue $ws "E42" "  -0.20%  "

Set-TextValue $ws "E43" "  -0.06%  "

Set-TextValue $ws "D44" "0.0489"
Set-TextValue $ws "E44" "  +2.73%  "

Set-TextValue $ws "E45" "  +6.41%  "

Set-TextValue $ws "D46" "2.64"
Set-TextValue $ws "E46" "  -6.59%  "

Set-TextValue $ws "D47" "9.61"
Set-TextValue $ws "E47" "  +9.51%  "

Set-TextValue $ws "E48" "  -7.34%  "

Set-TextValue $ws "D49" "3.05"
Set-TextValue $ws "E49" "  -0.34%  "

Set-TextValue $ws "D50" "3.35"
Set-TextValue $ws "E50" "  +1.10%  "

Set-TextValue $ws "D51" "144.47"
Set-TextValue $ws "E51" "  +0.62%  "
